$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Systems")

# 1) Insert a new row at row 5 (pushes "Caltech Radar I" and everything below down by one)
$ws.Rows("5:5").Insert()

# 2) Update row 4 ("Jena I") values that changed, and add the new B4 value.
#    Order matters for shared-string append order: 8.5 m, then 1 m diam,
#    then "Jena II" (row 5), then the new B4 frequency string.
$ws.Range("F4").Value = "8.5 m"
$ws.Range("G4").Value = "1 m diam"
$ws.Range("I4").Value = "0.4 K"

# 3) Populate the newly inserted row 5 ("Jena II") - mirrors Jena I's new values.
$ws.Range("A5").Value = "Jena II"
$ws.Range("C5").Value = 20
$ws.Range("D5").Value = "10 FPS"
$ws.Range("E5").Value = "0.5 K"
$ws.Range("F5").Value = "8.5 m"
$ws.Range("G5").Value = "1 m diam"
$ws.Range("H5").Value = "2cm"
$ws.Range("I5").Value = "0.4 K"

# 4) Now that row5 content exists, set the new B4 value (must be last of the new strings).
$ws.Range("B4").Value = "350 GHz / 80 GHz = 23 %"

# 5) Fix up hyperlinks: Insert() does not shift hyperlink ranges, so clear them all
#    and re-add at the correct (shifted) destinations.
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("J16"), "http://globenewswire.com/news-release/2011/07/19/451549/226833/en/Microsemi-Acquires-the-Technology-and-Related-Assets-From-Brijot-Imaging-Systems-Inc.html")
$ws.Hyperlinks.Add($ws.Range("J13"), "http://www.leidos.com/products/security/counterbomber")
$ws.Hyperlinks.Add($ws.Range("J14"), "http://www.arrowmid.com/Products/pdf/BIS-WDS_GEN2_Cut_Sheet.pdf")
$ws.Hyperlinks.Add($ws.Range("K9"), "http://iopscience.iop.org/1742-6596/400/5/052018/")
